$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155, shifting the existing rows 155-178 down
# to 156-179 (this also updates the sheet dimension automatically).
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new weekly data point,
# copying the static/categorical values from the (now shifted) row below
# and filling the updated numeric/date fields.
$ws.Range("A155").Value = 9
$ws.Range("B155").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C155").Value = "Metropolitana"
$ws.Range("D155").Value = 44505
$ws.Range("E155").Value = 13
$ws.Range("F155").Value = 300000001
$ws.Range("G155").Value = "Rabanito"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 7900
$ws.Range("K155").Value = 2500
$ws.Range("L155").Value = 3000
$ws.Range("M155").Value = 2747
$ws.Range("N155").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O155").Value = "Provincia de Chacabuco"
$ws.Range("P155").Value = 27
$ws.Range("Q155").Value = 100
$ws.Range("R155").Value = "Hortaliza"
